# Backup before dimension reduction
# Shift the "q" index labels in column A (rows 2-97) down by one:
#   A2: q1  -> q0
#   A3: q2  -> q1
#   ...
#   A97: q96 -> q95

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $newIndex = $row - 2
    $ws.Cells.Item($row, 1).Value = "q$newIndex"
}
